# [fix] id into excel
# Renumber the id_authorization column (A) on the active sheet so the
# sequence becomes contiguous starting at 1600 (row 2) through 1674
# (row 76), and move the active selection to K66 (scrolled so row 48
# is the top visible row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 76
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r + 1598
}

# Update the view: scroll to row 48 and select K66, matching the
# sheetView/selection state recorded in the saved workbook.
$ws.Range("K66").Select()
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 1
